# Swap the two embedded DrawingML colour themes (the presentation / slide-master
# theme and the notes-master theme traded places) and change the table style
# applied to the three tables in the deck to the "No Style, Table Grid" style.

$p = $ppt.ActivePresentation

# --- 1. Re-point the three tables at the new table style -------------------
$tableSlides = @(14, 15, 16)
foreach ($slideIdx in $tableSlides) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{C0548974-8186-437B-B2DC-E9D1EDA82F34}")
        }
    }
}

# --- 2. Swap the colour scheme used by the presentation theme --------------
# Before: the slide master's theme carried the "Red Violet" (Integral) colours
# and the notes master's theme carried the "Office" colours. After the edit
# the slide master's theme uses the "Office" colours (what used to live in
# the notes-master theme).
$slide = $p.Slides.Item(14)
$colors = $slide.ThemeColorScheme
$colors.Item(1).RGB  = 0          # dk1       000000
$colors.Item(2).RGB  = 16777215   # lt1       FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2       44546A
$colors.Item(4).RGB  = 15132391   # lt2       E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1   5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2   ED7D31
$colors.Item(7).RGB  = 10855845   # accent3   A5A5A5
$colors.Item(8).RGB  = 49407      # accent4   FFC000
$colors.Item(9).RGB  = 12874308   # accent5   4472C4
$colors.Item(10).RGB = 4697456    # accent6   70AD47
$colors.Item(11).RGB = 12673797   # hlink     0563C1
$colors.Item(12).RGB = 7491477    # folHlink  954F72
